# Performance results obtained for all java benchmarks in renaissance.
#
# Adds summary rows (AVERAGE / STDEV.P) below the existing data, highlights
# them with fills, widens columns A:B to fit content, and moves the active
# selection to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 21: column averages -----------------------------------------
$ws.Range("A21").Formula = "=AVERAGE(A2:A20)"
$ws.Range("B21").Formula = "=AVERAGE(B2:B20)"

# --- New row 22: population standard deviation ----------------------------
$ws.Range("A22").Formula = "=STDEV.P(A2:A20)"
$ws.Range("B22").Formula = "=STDEV.P(B2:B20)"

# --- Highlight the new summary rows ---------------------------------------
# Row 21 -> solid yellow fill
$ws.Range("A21:B21").Interior.Color = 65535
# Row 22 -> solid theme Accent2 fill
$ws.Range("A22:B22").Interior.Color = 65535
$ws.Range("A22:B22").Interior.ThemeColor = 6

# --- Widen columns A:B to fit the (now wider) content ---------------------
$ws.Columns("A:B").ColumnWidth = 11.1666666666667

# --- Move the active selection ---------------------------------------------
$ws.Range("D11").Select() | Out-Null
